$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.08045581118042906
$ws.Range("C2").Value = 0.8732618137668521
$ws.Range("D2").Value = 1.58169622500375
$ws.Range("E2").Value = 1.25765505008478
$ws.Range("F2").Value = 1.286073186806879

$ws.Range("B3").Value = 0.2468526762731521
$ws.Range("C3").Value = 1.012662856856466
$ws.Range("D3").Value = 1.831719206154512
$ws.Range("E3").Value = 1.353410213554823
$ws.Range("F3").Value = 1.362022745984608

$ws.Range("B4").Value = 0.2298317591110559
$ws.Range("C4").Value = 1.539245905603442
$ws.Range("D4").Value = 10.85920011347265
$ws.Range("E4").Value = 3.295330046212769
$ws.Range("F4").Value = 3.368485765113963
